# Fix figure/table placement on the "Content with Caption" slide layout
# (CustomLayout #8 of the slide master) by repositioning three shapes:
#   - the title placeholder
#   - the content placeholder (idx=1), which is also resized
#   - the caption/body placeholder (idx=2)
#
# PowerPoint COM exposes Left/Top/Width/Height in points (Single/float32),
# while OOXML stores EMUs (1 pt = 12700 EMU). The point values below are
# chosen so that, after the host's internal float32 round-trip back to
# EMU, they land exactly on the target EMU coordinates from the diff.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layout = $master.CustomLayouts.Item(8)

# Shape 1: title placeholder "標題 1" -> off x=671624 y=446690 (EMU)
$titleShape = $layout.Shapes.Item(1)
$titleShape.Left = 52.88377952755906
$titleShape.Top = 35.17244194488189

# Shape 2: content placeholder "內容版面配置區 2" (idx=1)
# -> off x=4332891 y=446690, ext cx=7564820 cy=5565228 (EMU)
$contentShape = $layout.Shapes.Item(2)
$contentShape.Left = 341.1725316850394
$contentShape.Top = 35.17244194488189
$contentShape.Width = 595.6551211102362
$contentShape.Height = 438.2069291338583

# Shape 3: text/body placeholder "文字版面配置區 3" (idx=2)
# -> off x=671624 y=1885465 (EMU)
$bodyShape = $layout.Shapes.Item(3)
$bodyShape.Left = 52.88377952755906
$bodyShape.Top = 148.46181502362202
